$wb = $excel.ActiveWorkbook

# --- Sheet "Info" ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Cells.Item(2, 1).Value = 640108574274.0112
$wsInfo.Cells.Item(2, 2).Value = 2.355999946594238

# --- Sheet "Activados" ---
# Grows from A1:B4 (3 data rows) to A1:B20 (19 data rows).
# Column A becomes 1 for every data row; column B becomes 0,20,40,...,360
$wsAct = $wb.Worksheets.Item("Activados")
for ($i = 0; $i -lt 19; $i++) {
    $row = $i + 2
    $wsAct.Cells.Item($row, 1).Value = 1
    $wsAct.Cells.Item($row, 2).Value = $i * 20
}

# --- Sheet "Operando" ---
# Column A becomes 1 for every data row (rows 2-366); column B (Tiempo) unchanged.
$wsOp = $wb.Worksheets.Item("Operando")
for ($row = 2; $row -le 366; $row++) {
    $wsOp.Cells.Item($row, 1).Value = 1
}

# --- Sheet "Contaminantes" ---
$wsCont = $wb.Worksheets.Item("Contaminantes")

$wsCont.Cells.Item(2, 2).Value = 449208244800.0004
$wsCont.Cells.Item(2, 3).Value = 16.66000000000001

$wsCont.Cells.Item(3, 2).Value = 13481640000.00001
$wsCont.Cells.Item(3, 3).Value = 0.5000000000000004

$wsCont.Cells.Item(4, 2).Value = 87091394399.99998
$wsCont.Cells.Item(4, 3).Value = 3.23

$wsCont.Cells.Item(5, 2).Value = 307074.010608
$wsCont.Cells.Item(5, 3).Value = 0.0000113886

$wsCont.Cells.Item(6, 2).Value = 90326988000.00008
$wsCont.Cells.Item(6, 3).Value = 3.350000000000003
